$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.494.01"
$ws.Range("E2").Value = "  -1.11%  "
$ws.Range("D3").Value = "3.984.56"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "541.52"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +4.15%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "149.29"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.28%  "
$ws.Range("D7").Value = "3.976.27"
$ws.Range("E7").Value = "  -1.44%  "
$ws.Range("E8").Value = "  -4.96%  "
$ws.Range("E9").Value = "  +0.00%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.746"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -4.01%  "
$ws.Range("E11").Value = "  -4.23%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "56.57"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +18.52%  "
$ws.Range("E13").Value = "  -2.54%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "10.77"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "4.616.10"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").Value = "3.977.75"
$ws.Range("E16").Value = "  -2.31%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "14.00"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.82%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "20.57"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -2.87%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.131"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.26%  "
$ws.Range("E20").Value = "  -2.83%  "
$ws.Range("D21").Value = "71.368.98"
$ws.Range("E21").Value = "  -1.05%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "429.21"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -3.25%  "
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "97.62"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -6.92%  "
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "3.60"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +1.32%  "
$ws.Range("E25").Value = "  +5.31%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "14.55"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -2.19%  "
$ws.Range("E27").Value = "  +0.44%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "10.78"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -2.13%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "3.76"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +16.15%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.92"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.30%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "36.70"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("E32").Value = "  +12.56%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "51.56"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +20.88%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.132"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +0.99%  "
$ws.Range("B35").Value = "Cosmos"
$ws.Range("C35").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "13.43"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("B36").Value = "Bittensor"
$ws.Range("C36").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "685.92"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.62%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "65.58"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -1.93%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.440"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.81%  "
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "0.0₃0823"
$ws.Range("E40").Value = "  -4.88%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "3.41"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.59%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.997"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -0.12%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +0.31%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.28"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  -2.73%  "
$ws.Range("E46").Value = "  -5.42%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.51%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "9.73"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +6.32%  "
$ws.Range("E49").Value = "  -5.50%  "
$ws.Range("E50").Value = "  -1.75%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.000272"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +2.39%  "
